$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 14:05"

# Row 6 - Brasil
$ws.Range("B6").Value = 312074
$ws.Range("C6").Value = 1153
$ws.Range("E6").Value = 166002
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 20112

# Row 27 - Suecia
$ws.Range("B27").Value = 32809
$ws.Range("C27").Value = 637
$ws.Range("E27").Value = 23913
$ws.Range("G27").Value = 54
$ws.Range("H27").Value = 3925

# Row 28 - Suiza
$ws.Range("B28").Value = 30707
$ws.Range("C28").Value = 13
$ws.Range("E28").Value = 904
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 1903

# Row 87 - becomes "Consejo Danes para los Refugiados"
$ws.Range("A87").Value = "Consejo Danes para los Refugiados"
$ws.Range("B87").Value = 1945
$ws.Range("C87").Value = 110
$ws.Range("D87").Value = 312
$ws.Range("E87").Value = 1570
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 63

# Row 88 - Republica de Macedonia (name unchanged, data updated)
$ws.Range("B88").Value = 1921
$ws.Range("C88").Value = 23
$ws.Range("D88").Value = 1387
$ws.Range("E88").Value = 422
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 112

# Row 89 - becomes "Cuba"
$ws.Range("A89").Value = "Cuba"
$ws.Range("B89").Value = 1908
$ws.Range("D89").Value = 1603
$ws.Range("E89").Value = 225
$ws.Range("H89").Value = 80

# Row 104 - becomes "Libano"
$ws.Range("A104").Value = "Libano"
$ws.Range("B104").Value = 1086
$ws.Range("C104").Value = 62
$ws.Range("D104").Value = 663
$ws.Range("E104").Value = 397
$ws.Range("H104").Value = 26

# Row 105 - becomes "Hong Kong"
$ws.Range("A105").Value = "Hong Kong"
$ws.Range("B105").Value = 1066
$ws.Range("C105").Value = 2
$ws.Range("D105").Value = 1029
$ws.Range("E105").Value = 33
$ws.Range("H105").Value = 4

# Row 106 - becomes "Sri Lanka"
$ws.Range("A106").Value = "Sri Lanka"
$ws.Range("B106").Value = 1055
$ws.Range("C106").Value = 7
$ws.Range("D106").Value = 620
$ws.Range("E106").Value = 426
$ws.Range("H106").Value = 9

# Row 107 - becomes "Tunez"
$ws.Range("A107").Value = "Tunez"
$ws.Range("B107").Value = 1046
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 883
$ws.Range("E107").Value = 116
$ws.Range("H107").Value = 47

# Row 108 - becomes "Letonia"
$ws.Range("A108").Value = "Letonia"
$ws.Range("B108").Value = 1030
$ws.Range("C108").Value = 5
$ws.Range("D108").Value = 712
$ws.Range("E108").Value = 296
$ws.Range("H108").Value = 22

# Row 113 - Guinea Ecuatorial (name unchanged, data updated)
$ws.Range("D113").Value = 165
$ws.Range("E113").Value = 728

# Row 194 - becomes "Namibia"
$ws.Range("A194").Value = "Namibia"
$ws.Range("C194").Value = 1

# Row 195 - becomes "Laos"
$ws.Range("A195").Value = "Laos"
$ws.Range("C195").Value = 0

# Row 198 - becomes "Santa Lucia"
$ws.Range("A198").Value = "Santa Lucia"
$ws.Range("D198").Value = 18
$ws.Range("H198").Value = 0

# Row 200 - becomes "Belice"
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

# Row 209 - becomes "Montserrat"
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

# Row 211 - becomes "Groenlandia"
$ws.Range("A211").Value = "Groenlandia"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
